# ---------------------------------------------------------------------------
# Testdata.xlsx update:
#   - data1 (sheet2): add ALERT_MESSAGE/RUN_STATUS values in D2:E6
#   - data2 (sheet3): add the two blank-alert columns + RUN_STATUS in D2:F6
#   - data3 (new sheet, copied tab layout of data1/data2): TEST_CASE_ID,
#     EMAIL_ID (with mailto hyperlinks), USERNAME, PASSWORD, RUN_STATUS
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. data1 : fill in the new ALERT_MESSAGE / RUN_STATUS columns (D, E)
# ---------------------------------------------------------------------------
$data1 = $wb.Worksheets.Item("data1")

$data1.Range("D2").Value = "User or Password is not valid"
$data1.Range("D3").Value = "User or Password is not valid"
$data1.Range("D4").Value = "User or Password is not valid"
$data1.Range("D5").Value = "User or Password is not valid"
$data1.Range("D6").Value = "User or Password is not valid"

$data1.Range("E2").Value = "PASSED"
$data1.Range("E3").Value = "PASSED"
$data1.Range("E4").Value = "PASSED"
$data1.Range("E5").Value = "PASSED"
$data1.Range("E6").Value = "PASSED"

$data1.Range("E2:E6").Select()

# ---------------------------------------------------------------------------
# 2. data2 : fill in the new blank-alert / RUN_STATUS columns (D, E, F)
# ---------------------------------------------------------------------------
$data2 = $wb.Worksheets.Item("data2")

$data2.Range("D2").Value = "User-ID must not be blank"
$data2.Range("D3").Value = "User-ID must not be blank"
$data2.Range("D4").Value = "User-ID must not be blank"
$data2.Range("D5").Value = "User-ID must not be blank"
$data2.Range("D6").Value = "User-ID must not be blank"

$data2.Range("E2").Value = "Password must not be blank"
$data2.Range("E3").Value = "Password must not be blank"
$data2.Range("E4").Value = "Password must not be blank"
$data2.Range("E5").Value = "Password must not be blank"
$data2.Range("E6").Value = "Password must not be blank"

$data2.Range("F2").Value = "PASSED"
$data2.Range("F3").Value = "PASSED"
$data2.Range("F4").Value = "PASSED"
$data2.Range("F5").Value = "PASSED"
$data2.Range("F6").Value = "PASSED"

$data2.Range("E13").Select()

# ---------------------------------------------------------------------------
# 3. data3 : brand-new worksheet, appended after data2
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$data3 = $wb.Worksheets.Add($null, $lastSheet)
$data3.Name = "data3"

# column widths (roughly matching the author's layout)
$data3.Columns.Item(1).ColumnWidth = 17.45
$data3.Columns.Item(2).ColumnWidth = 27.02
$data3.Columns.Item(3).ColumnWidth = 27.02
$data3.Columns.Item(4).ColumnWidth = 20.17
$data3.Columns.Item(5).ColumnWidth = 23.31

# header row - reuse the same formatting as the other data sheets
$data1.Range("A1:E1").Copy()
$data3.Range("A1:E1").PasteSpecial(-4122)

$data3.Range("A1").Value = "TEST_CASE_ID"
$data3.Range("B1").Value = "EMAIL_ID"
$data3.Range("C1").Value = "USERNAME"
$data3.Range("D1").Value = "PASSWORD"
$data3.Range("E1").Value = "RUN_STATUS"

# column A (TEST_CASE_ID) - reuse the bordered body-cell style
$data1.Range("A2:A6").Copy()
$data3.Range("A2:A6").PasteSpecial(-4122)

$data3.Range("A2").Value = "TC_21"
$data3.Range("A3").Value = "TC_22"
$data3.Range("A4").Value = "TC_23"
$data3.Range("A5").Value = "TC_24"
$data3.Range("A6").Value = "TC_25"

# column B (EMAIL_ID) - same bordered style, then turned into hyperlinks
$data1.Range("A2:A6").Copy()
$data3.Range("B2:B6").PasteSpecial(-4122)

$data3.Range("B2").Value = "kaushik.0407@gmail.com"
$data3.Range("B3").Value = "kaushik.0407@gmail.com"
$data3.Range("B4").Value = "kaushik.0407@gmail.com"
$data3.Range("B5").Value = "kaushik.0407@gmail.com"
$data3.Range("B6").Value = "kaushik.0407@gmail.com"

$data3.Hyperlinks.Add($data3.Range("B2"), "mailto:kaushik.0407@gmail.com")
$data3.Hyperlinks.Add($data3.Range("B3:B6"), "mailto:kaushik.0407@gmail.com", "", "", "kaushik.0407@gmail.com")

$data3.Range("B3").Value = "kaushik.0407@yahoo.com"
$data3.Hyperlinks.Add($data3.Range("B3"), "mailto:kaushik.0407@yahoo.com")

$data3.Range("B4").Value = "kaushik.0407@rediffmail.com"
$data3.Hyperlinks.Add($data3.Range("B4"), "mailto:kaushik.0407@rediffmail.com")

$data3.Range("B5").Value = "kaushik.0407@hotmail.com"
$data3.Hyperlinks.Add($data3.Range("B5"), "mailto:kaushik.0407@hotmail.com")

$data3.Range("B6").Value = "kaushik.0407@apple.com"
$data3.Hyperlinks.Add($data3.Range("B6"), "mailto:kaushik.0407@apple.com")

# columns C, D, E (USERNAME / PASSWORD / RUN_STATUS) - plain, unstyled cells
$data3.Range("C2").Value = "mngr299505"
$data3.Range("D2").Value = "nYdAvun"
$data3.Range("E2").Value = "PASSED"

$data3.Range("C3").Value = "mngr299523"
$data3.Range("D3").Value = "AhujYtu"
$data3.Range("E3").Value = "PASSED"

$data3.Range("C4").Value = "mngr299525"
$data3.Range("D4").Value = "vEvUmAr"
$data3.Range("E4").Value = "PASSED"

$data3.Range("C5").Value = "mngr299524"
$data3.Range("D5").Value = "ArEbUmu"
$data3.Range("E5").Value = "PASSED"

$data3.Range("C6").Value = "mngr299515"
$data3.Range("D6").Value = "durugEj"
$data3.Range("E6").Value = "PASSED"

# final selection / active sheet - matches the author leaving the cursor on
# the brand-new data3 sheet with C2:E6 selected
$data3.Range("C2:E6").Select()
